$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date values for rows 2-8 from 45185 (2023-09-16)
# to 45204 (2023-10-05), matching the automatic update of the logging report.
$ws.Range("C2:C8").Value = 45204
